$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Mahmoud"
$ws.Range("B2").Value = "Galal"

# Phone number has a leading zero; force text storage so it is not
# coerced into a number (which would drop the leading "0"), then drop
# the formatting back to the sheet's normal style.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "01096001734"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "Full Stack Developer"
$ws.Range("E2").Value = "admin@admin.com"
$ws.Range("F2").Value = "awdjwadioajd , dwajwaopidjawodpawjd, awdawdjawopid"
